# Delete kingdom & taxonRank from data template
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("checklist")

# Remove columns B (kingdom) and C (taxon_rank), shifting remaining columns left
$ws.Range("B:C").Delete()

# Update the active selection to match the post-edit workbook state
$ws.Range("D22").Select()
